$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several updated Price values are plain decimal numbers (e.g. "1.000",
# "299.75"). Pre-format the whole Price column as Text so Excel keeps the
# literal digits (incl. trailing zeros) instead of auto-converting them to
# a number when the new value is written below.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "23.505.75"
$ws.Range("E2").Value = "  -1.01%  "

# Row 3
$ws.Range("D3").Value = "1.649.71"
$ws.Range("E3").Value = "  -0.25%  "

# Row 4
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.30%  "

# Row 5
$ws.Range("E5").Value = "  +0.21%  "

# Row 6
$ws.Range("D6").Value = "299.75"
$ws.Range("E6").Value = "  -1.68%  "

# Row 7
$ws.Range("D7").Value = "0.3796"
$ws.Range("E7").Value = "  -0.71%  "

# Row 8
$ws.Range("D8").Value = "0.3563"
$ws.Range("E8").Value = "  -1.65%  "

# Row 9
$ws.Range("D9").Value = "50.12"
$ws.Range("E9").Value = "  -3.70%  "

# Row 10
$ws.Range("D10").Value = "0.08106"
$ws.Range("E10").Value = "  -1.56%  "

# Row 11
$ws.Range("D11").Value = "1.222"
$ws.Range("E11").Value = "  -2.32%  "

# Row 12
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  +0.29%  "

# Row 13
$ws.Range("D13").Value = "22.06"
$ws.Range("E13").Value = "  -2.38%  "

# Row 14
$ws.Range("D14").Value = "6.408"
$ws.Range("E14").Value = "  -1.97%  "

# Row 15
$ws.Range("D15").Value = "7.389"
$ws.Range("E15").Value = "  -0.04%  "

# Row 16
$ws.Range("D16").Value = "0.00001196"
$ws.Range("E16").Value = "  -2.91%  "

# Row 17
$ws.Range("D17").Value = "1.654.11"
$ws.Range("E17").Value = "  -0.22%  "

# Row 18
$ws.Range("D18").Value = "97.41"
$ws.Range("E18").Value = "  +0.56%  "

# Row 19
$ws.Range("D19").Value = "0.06966"
$ws.Range("E19").Value = "  +0.02%  "

# Row 20
$ws.Range("D20").Value = "6.775"
$ws.Range("E20").Value = "  -0.35%  "

# Row 21
$ws.Range("D21").Value = "17.24"
$ws.Range("E21").Value = "  -2.57%  "

# Row 22
$ws.Range("E22").Value = "  +0.12%  "

# Row 23
$ws.Range("D23").Value = "12.44"
$ws.Range("E23").Value = "  -1.47%  "

# Row 24
$ws.Range("D24").Value = "23.529.62"
$ws.Range("E24").Value = "  -0.88%  "

# Row 25
$ws.Range("D25").Value = "2.502"
$ws.Range("E25").Value = "  -0.92%  "

# Row 26
$ws.Range("D26").Value = "2.912"
$ws.Range("E26").Value = "  -5.90%  "

# Row 27
$ws.Range("E27").Value = "  -1.83%  "

# Row 28
$ws.Range("D28").Value = "153.69"

# Row 29
$ws.Range("D29").Value = "5.200"
$ws.Range("E29").Value = "  -0.12%  "

# Row 30
$ws.Range("D30").Value = "132.75"
$ws.Range("E30").Value = "  -1.50%  "

# Row 31
$ws.Range("D31").Value = "1.837.34"
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("D32").Value = "6.929"
$ws.Range("E32").Value = "  +0.60%  "

# Row 33
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "2.120"
$ws.Range("E33").Value = "  +0.82%  "

# Row 34
$ws.Range("B34").Value = "FraxShare"
$ws.Range("C34").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D34").Value = "11.89"
$ws.Range("E34").Value = "  +1.19%  "

# Row 35
$ws.Range("E35").Value = "  -6.75%  "

# Row 36
$ws.Range("D36").Value = "0.02730"
$ws.Range("E36").Value = "  -2.94%  "

# Row 37
$ws.Range("D37").Value = "0.08743"
$ws.Range("E37").Value = "  -0.87%  "

# Row 38
$ws.Range("D38").Value = "0.2440"
$ws.Range("E38").Value = "  -3.13%  "

# Row 39
$ws.Range("D39").Value = "5.962"
$ws.Range("E39").Value = "  -2.47%  "

# Row 40
$ws.Range("D40").Value = "13.20"
$ws.Range("E40").Value = "  +3.30%  "

# Row 41
$ws.Range("D41").Value = "0.06812"
$ws.Range("E41").Value = "  -3.49%  "

# Row 42
$ws.Range("D42").Value = "0.6905"
$ws.Range("E42").Value = "  -2.48%  "

# Row 43
$ws.Range("D43").Value = "1.318"
$ws.Range("E43").Value = "  -1.67%  "

# Row 44
$ws.Range("D44").Value = "15.52"
$ws.Range("E44").Value = "  -2.58%  "

# Row 45
$ws.Range("E45").Value = "  +0.19%  "

# Row 46
$ws.Range("D46").Value = "0.6413"
$ws.Range("E46").Value = "  -1.66%  "

# Row 47
$ws.Range("D47").Value = "2.262"
$ws.Range("E47").Value = "  -3.39%  "

# Row 48
$ws.Range("D48").Value = "3.922"
$ws.Range("E48").Value = "  -1.67%  "

# Row 49
$ws.Range("D49").Value = "0.07732"
$ws.Range("E49").Value = "  -3.24%  "

# Row 50
$ws.Range("D50").Value = "127.64"
$ws.Range("E50").Value = "  -0.58%  "

# Row 51
$ws.Range("D51").Value = "1.149"
$ws.Range("E51").Value = "  -3.75%  "
